$wb = $excel.ActiveWorkbook

# --- Sheet renames (workbook.xml) ---
$wsWeeknr = $wb.Worksheets.Item(1)
$wsWeeknr.Name = "weeknr 48"

$wsWeek = $wb.Worksheets.Item(2)
$wsWeek.Name = "week 49"

$wsTotaal = $wb.Worksheets.Item("Totaal")

# --- "Totaal" sheet ---
$wsTotaal.Range("B2").Value = "PyramidPanic"
$wsTotaal.Range("A7").Value = 48
$wsTotaal.Range("B7").Formula = "='weeknr 48'!G26"

# Update the active selection to B3 as in the diff
$wsTotaal.Range("B3").Select()

# --- "weeknr 48" sheet (was "weeknr 45") ---

# Row 10
$wsWeeknr.Range("C10").Value = 0.4375
$wsWeeknr.Range("D10").Value = 0.4513888888888889
$wsWeeknr.Range("F10").Value = "Alle essets geordend en github + logboek bijgewerkt."

# Row 11
$wsWeeknr.Range("C11").Value = 0.4548611111111111
$wsWeeknr.Range("D11").Value = 0.46527777777777773
$wsWeeknr.Range("F11").Value = "Spritebatch toegevoegd voor background."

# Row 12
$wsWeeknr.Range("C12").Value = 0.46875
$wsWeeknr.Range("D12").Value = 0.4861111111111111
$wsWeeknr.Range("F12").Value = "IsMouseVisble = true gemaakt en escape toegevoegd."

# Update the active selection to F12 as in the diff, and make this the active sheet
$wsWeeknr.Activate()
$wsWeeknr.Range("F12").Select()

$wb.Save()
